$d = $word.ActiveDocument

# Each of the 5 answer bullets below needs "Design: " prepended to the
# start of their text.
# (Table.Cells indexing is unreliable for multi-paragraph cells in this
# runtime, so the target paragraphs are addressed directly through
# Document.Paragraphs, identified by their current leading text.)

$targets = @(
    "systeemit toimivat suunnitellun mukaisesti.",
    "Projektissa oli mielestäni realistisempi aikataulu kuin aikasemmissa",
    "Toimi melko hyvin",
    "Toimi",
    "Aikataulu"
)

$cr = [char]13   # paragraph mark
$vt = [char]11   # manual line break (w:br)
$bel = [char]7   # end-of-cell mark

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    foreach ($target in $targets) {
        if ($text.StartsWith($target)) {
            $rest = $text.Substring($target.Length)
            $nextChar = ""
            if ($rest.Length -gt 0) {
                $nextChar = $rest.Substring(0,1)
            }
            $isExactMatch = ($rest.Length -eq 0) -or ($nextChar -eq $cr) -or ($nextChar -eq $vt) -or ($nextChar -eq $bel)
            if ($isExactMatch) {
                $p.Range.InsertBefore("Design: ")
                break
            }
        }
    }
}
